$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New filter columns (K, L, M) added to the flight-search header/row data.
# Values are written in this exact order so the shared-string table is
# built up the same way the original authoring session produced it.
$ws.Range("K1").Value = "Stopfilter"
$ws.Range("K2").Value = "Non-Stop"
$ws.Range("L1").Value = "Airlinesfilter"
$ws.Range("M1").Value = "Departuretimefilter"
$ws.Range("M2").Value = "EARLY_MORNING"
$ws.Range("L2").Value = "6E "

# New column M needs an explicit width like the other data columns.
$ws.Columns.Item(13).ColumnWidth = 21.83

# Selection / scroll position moved as the user worked on the new columns.
[void]$ws.Range("K5").Select()
$excel.ActiveWindow.ScrollColumn = 5

# Page now explicitly set to portrait orientation for printing.
$ws.PageSetup.Orientation = 1
